$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# SCD0205 - Validasi Field report PHR pada searching portal
# Two new columns ("TEXT4", "TEXT5") are inserted right before the
# existing "FILE1" column (which was column O). Since a real
# EntireColumn.Insert shifts the whole sheet (including the empty
# formatting rows below the data, which must stay put here), the
# columns are instead rebuilt cell-by-cell:
#   - the old "FILE1" column (O) content/format is moved to column Q
#   - column O is reused for the new "TEXT4" header
#   - column P becomes the new "TEXT5" header
# ---------------------------------------------------------------------

# Remember the width of the old FILE1 column (O) and of an existing
# TEXT column (L) before anything is changed.
$file1Width = $ws.Columns("O").ColumnWidth
$textColWidth = $ws.Columns("L").ColumnWidth

# --- Move the FILE1 data column (O) to its new home in column Q ------

# Row 2 (empty data cell, keep its number format / vertical alignment)
$ws.Range("O2").Copy()
$ws.Range("Q2").PasteSpecial(-4122)  # xlPasteFormats

# Row 3 ("Februari" file value, keep format + bring the value along)
$ws.Range("O3").Copy()
$ws.Range("Q3").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("Q3").Value2 = $ws.Range("O3").Value2

$ws.Application.CutCopyMode = 0

# Row 1 header: FILE1 moves from O1 to Q1
$ws.Range("Q1").Value2 = $ws.Range("O1").Value2

# Clear the now-vacated old cells in column O (rows 2-3)
$ws.Range("O2:O3").Clear()

# --- Fill in the two new header cells ---------------------------------
$ws.Range("O1").Value2 = "TEXT4"
$ws.Range("P1").Value2 = "TEXT5"

# --- Column widths ------------------------------------------------------
# New TEXT4 / TEXT5 columns get the same width as the other TEXT columns,
# and the relocated FILE1 column (Q) keeps its original width.
$ws.Columns("O").ColumnWidth = $textColWidth
$ws.Columns("P").ColumnWidth = $textColWidth
$ws.Columns("Q").ColumnWidth = $file1Width

# --- Selection ------------------------------------------------------
$ws.Range("Q2").Select()
